$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = (Get-Date -Year 2016 -Month 8 -Day 30 -Hour 21 -Minute 29 -Second 20)
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 53
$ws.Range("D6").Value = 44
$ws.Range("E6").Value = 20
$ws.Range("F6").Value = 80
$ws.Range("G6").Value = 16398
$ws.Range("H6").Value = 11622
$ws.Range("I6").Value = 1786
$ws.Range("J6").Value = 267
$ws.Range("K6").Value = 220
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 4
$ws.Range("N6").Value = "Noun"
